# Use case workbook update: add "encoding" sheet with hierarchy encoding
# reference table, and restyle the GA4 MERCH example table on Folha1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Add the new "encoding" sheet right after Folha1 and make it active.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "encoding"

# ---------------------------------------------------------------------
# 2) Populate the "encoding" sheet.
# ---------------------------------------------------------------------
$ws2.Range("I11").Value = "hierarchy enconding"
$ws2.Range("I12").Value = "A-Z order?"

$ws2.Range("I15").Value = "cat1"
$ws2.Range("J15").Value = "cat2"
$ws2.Range("K15").Value = "brand"
$ws2.Range("L15").Value = "cat encoded"
$ws2.Range("M15").Value = "cat 2 encoded"
$ws2.Range("N15").Value = "brand encoded"

$ws2.Range("I16").Value = "electronics"
$ws2.Range("J16").Value = "tv"
$ws2.Range("K16").Value = "lg"
$ws2.Range("L16").Value = 1
$ws2.Range("M16").Value = 1
$ws2.Range("N16").Value = 1

$ws2.Range("I17").Value = "electronics"
$ws2.Range("J17").Value = "smartphone"
$ws2.Range("K17").Value = "apple"
$ws2.Range("L17").Value = 1
$ws2.Range("M17").Value = 2
$ws2.Range("N17").Value = 2

$ws2.Range("I18").Value = "home appliances"
$ws2.Range("J18").Value = "fridge"
$ws2.Range("K18").Value = "whirlpool"
$ws2.Range("L18").Value = 2
$ws2.Range("M18").Value = 3
$ws2.Range("N18").Value = 3

$ws2.Range("I19").Value = "home appliances"
$ws2.Range("J19").Value = "fridge"
$ws2.Range("K19").Value = "lg"
$ws2.Range("L19").Value = 2
$ws2.Range("M19").Value = 3
$ws2.Range("N19").Value = 1

$ws2.Range("I20").Value = "apparell"
$ws2.Range("J20").Value = "tshirts"
$ws2.Range("K20").Value = "nike"
$ws2.Range("L20").Value = 3
$ws2.Range("M20").Value = 4
$ws2.Range("N20").Value = 4

# Row 15 header cells (L:N) get wrap-text formatting, and the whole header
# row is taller to fit the wrapped text.
$ws2.Range("L15:N15").WrapText = $true
$ws2.Rows.Item(15).RowHeight = 30

# Column I is a bit wider to fit "hierarchy enconding" / "home appliances"
# (stored column width of 15.57 characters).
$ws2.Columns.Item(9).ColumnWidth = 14.74

# ---------------------------------------------------------------------
# 3) View/selection tweaks to match the authoring session: "encoding"
#    becomes the active/selected sheet, zoomed to 270%, with M15 selected.
# ---------------------------------------------------------------------
$ws2.Activate()
$excel.ActiveWindow.Zoom = 270
$ws2.Range("M15").Select()
